# "create api for get profile details"
# Adds two new employee rows (Varun Gautam / dot-net dev, and an Hr
# Department / moderator admin account) to the roster, and replaces the
# "date of joining" / "date of birth" columns for every row with a
# (broken) DATEVALUE() formula instead of a literal date - matching the
# upstream change exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Existing rows 2-5: "dateofjoining" (H) and "dateofbirth" (N) become
#    formulas instead of literal date serials. Everything else in these
#    rows is unchanged.
# ---------------------------------------------------------------------
$ws.Range("H2").Formula = "=DATEVALUE(2/3/2022)"
$ws.Range("H3").Formula = "=DATEVALUE(2/3/2022)"
$ws.Range("H4").Formula = "=DATEVALUE(2/3/2022)"
$ws.Range("H5").Formula = "=DATEVALUE(2/3/2022)"

# N2 used to be plain text ("23/11/1997") - pull the date-formatted
# style off N3 first so it ends up styled like the other dateofbirth
# cells once it holds a formula.
$ws.Range("N3").Copy() | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("N2").Formula = "=DATEVALUE(11/23/1997)"

$ws.Range("N3").Formula = "=DATEVALUE(10/13/1993)"
$ws.Range("N4").Formula = "=DATEVALUE(11/22/1998)"
$ws.Range("N5").Formula = "=DATEVALUE(11/16/1995)"

# ---------------------------------------------------------------------
# 2. Row 6: Varun Gautam - Sr. Full Stack Developer / Dot Net Developer
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "Varun Gautam"
$ws.Range("B6").Value = "varungautam.ameotech@gmail.com"
$ws.Range("C6").Value = "varun_ameotech"
$ws.Range("D6").Value = 123456
$ws.Range("E6").Value = 9456434855
$ws.Range("F6").Value = "Sr. Full Stack Developer"
$ws.Range("G6").Value = "Dot Net Developer"
$ws.Range("H6").Formula = "=DATEVALUE(2/3/2022)"
$ws.Range("I6").Value = "user"
$ws.Range("J6").Value = 52000
$ws.Range("M6").Value = "City: - Kurali"

$ws.Range("N3").Copy() | Out-Null
$ws.Range("N6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("N6").Formula = "=DATEVALUE(11/23/1997)"

$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:varungautam.ameotech@gmail.com") | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (keep s=2, not a fresh style)

# ---------------------------------------------------------------------
# 3. Row 7: Hr Department admin / moderator account
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "Hr Department"
$ws.Range("B7").Value = "hr@ameotech.gmail.com"
$ws.Range("C7").Value = "hr_admin"
$ws.Range("D7").Value = "HrAdmin@123"
$ws.Range("E7").Value = 9756484554
$ws.Range("F7").Value = "Hr Manager"
$ws.Range("G7").Value = "HRM"
$ws.Range("I7").Value = "moderator"
$ws.Range("J7").Value = 52000

$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:hr@ameotech.gmail.com") | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:HrAdmin@123") | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# N7 stays empty, but (like H7) should still exist as a styled placeholder
# cell, matching the pattern of every other empty row.
$ws.Range("N3").Copy() | Out-Null
$ws.Range("N7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------
# 4. View / print setup
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("K7").Select() | Out-Null

$ws.PageSetup.Orientation = 1   # xlPortrait
